$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), O (Origen), P (Precio $/Kg)
# derived from matching the canonical OOXML diff - the edit reshuffles these
# seven columns amongst the 25 data rows (a pure permutation; A,B,C,E,F,G,H,I,N,Q,R
# are identical across all rows so they are unaffected).
$rows = @(
    @{ Row = 2; D = 44202; J = 50; K = 8000; L = 9000; M = 8400; O = "Región del Maule"; P = 140 }
    @{ Row = 3; D = 44208; J = 100; K = 7000; L = 8000; M = 7350; O = "Región del Maule"; P = 122 }
    @{ Row = 4; D = 44617; J = 100; K = 10000; L = 11000; M = 10500; O = "Región Metropolitana"; P = 175 }
    @{ Row = 5; D = 44204; J = 45; K = 9500; L = 10000; M = 9722; O = "Región del Maule"; P = 162 }
    @{ Row = 6; D = 44264; J = 43; K = 8500; L = 9000; M = 8709; O = "Región del Maule"; P = 145 }
    @{ Row = 7; D = 44627; J = 60; K = 9000; L = 9500; M = 9250; O = "Región Metropolitana"; P = 154 }
    @{ Row = 8; D = 44210; J = 60; K = 8000; L = 9000; M = 8417; O = "Región de Arica y Parinacota"; P = 140 }
    @{ Row = 9; D = 44253; J = 95; K = 9500; L = 10000; M = 9658; O = "Región del Maule"; P = 161 }
    @{ Row = 10; D = 44610; J = 60; K = 11000; L = 12000; M = 11500; O = "Región Metropolitana"; P = 192 }
    @{ Row = 11; D = 44615; J = 100; K = 11000; L = 12000; M = 11500; O = "Región Metropolitana"; P = 192 }
    @{ Row = 12; D = 44624; J = 60; K = 10000; L = 11000; M = 10500; O = "Región Metropolitana"; P = 175 }
    @{ Row = 13; D = 44630; J = 60; K = 9000; L = 9500; M = 9250; O = "Región Metropolitana"; P = 154 }
    @{ Row = 14; D = 44218; J = 65; K = 9000; L = 10000; M = 9615; O = "Región del Maule"; P = 160 }
    @{ Row = 15; D = 44608; J = 100; K = 12000; L = 13000; M = 12500; O = "Región Metropolitana"; P = 208 }
    @{ Row = 16; D = 44160; J = 90; K = 7500; L = 8000; M = 7667; O = "Región de Arica y Parinacota"; P = 128 }
    @{ Row = 17; D = 44216; J = 55; K = 9500; L = 10000; M = 9773; O = "Región del Maule"; P = 163 }
    @{ Row = 18; D = 44159; J = 35; K = 7500; L = 8000; M = 7714; O = "Región de Arica y Parinacota"; P = 129 }
    @{ Row = 19; D = 44271; J = 55; K = 9000; L = 9500; M = 9227; O = "Región del Maule"; P = 154 }
    @{ Row = 20; D = 44162; J = 43; K = 8000; L = 8500; M = 8209; O = "Región de Arica y Parinacota"; P = 137 }
    @{ Row = 21; D = 44596; J = 100; K = 12000; L = 13000; M = 12500; O = "Región de Arica y Parinacota"; P = 208 }
    @{ Row = 22; D = 44600; J = 60; K = 12000; L = 13000; M = 12500; O = "Región de Arica y Parinacota"; P = 208 }
    @{ Row = 23; D = 44266; J = 60; K = 9000; L = 9500; M = 9208; O = "Región del Maule"; P = 153 }
    @{ Row = 24; D = 44224; J = 80; K = 8500; L = 9000; M = 8719; O = "Región del Maule"; P = 145 }
    @{ Row = 25; D = 44594; J = 80; K = 12000; L = 13000; M = 12500; O = "Región de Arica y Parinacota"; P = 208 }
    @{ Row = 26; D = 44259; J = 70; K = 9000; L = 9500; M = 9214; O = "Región del Maule"; P = 154 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D    # D: Fecha
    $ws.Cells.Item($r.Row, 10).Value = $r.J   # J: Volumen
    $ws.Cells.Item($r.Row, 11).Value = $r.K   # K: Precio minimo
    $ws.Cells.Item($r.Row, 12).Value = $r.L   # L: Precio maximo
    $ws.Cells.Item($r.Row, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r.Row, 15).Value = $r.O   # O: Origen
    $ws.Cells.Item($r.Row, 16).Value = $r.P   # P: Precio $/Kg
}
